$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 127, pushing existing rows 127-187 down to 128-188.
$ws.Rows.Item(127).Insert()

# Populate the newly inserted row 127 with the new weekly price record.
$ws.Cells.Item(127, 1).Value() = 11
$ws.Cells.Item(127, 2).Value() = "Vega Monumental Concepción"
$ws.Cells.Item(127, 3).Value() = "Bíobío"
$ws.Cells.Item(127, 4).Value() = 45009
$ws.Cells.Item(127, 5).Value() = 8
$ws.Cells.Item(127, 6).Value() = "Fruta"
$ws.Cells.Item(127, 7).Value() = 100102
$ws.Cells.Item(127, 8).Value() = "Cítricos"
$ws.Cells.Item(127, 9).Value() = 100102004
$ws.Cells.Item(127, 10).Value() = "Mandarina"
$ws.Cells.Item(127, 11).Value() = "Murcott"
$ws.Cells.Item(127, 12).Value() = "Primera"
$ws.Cells.Item(127, 13).Value() = 330
$ws.Cells.Item(127, 14).Value() = 8000
$ws.Cells.Item(127, 15).Value() = 9000
$ws.Cells.Item(127, 16).Value() = 8545
$ws.Cells.Item(127, 17).Value() = "`$/caja 15 kilos granel"
$ws.Cells.Item(127, 18).Value() = "Región de O'Higgins"
$ws.Cells.Item(127, 19).Value() = 570
$ws.Cells.Item(127, 20).Value() = 15
